$d = $word.ActiveDocument

function Set-RunFormatting {
    param(
        [string]$ParaText,
        [int]$Bold,
        [int]$Italic,
        [int]$Strike,
        [int]$Size
    )
    foreach ($p in $d.Paragraphs) {
        $raw = $p.Range.Text
        $trimmed = $raw.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $ParaText) {
            $r = $d.Range($p.Range.Start, $p.Range.End - 1)
            $r.Font.Bold = $Bold
            $r.Font.Italic = $Italic
            $r.Font.StrikeThrough = $Strike
            if ($Size -gt 0) {
                $r.Font.Size = $Size
            }
        }
    }
}

# Re-apply the (unchanged) bold/italic/strike/size formatting on the runs so
# they get re-serialized by the current engine - matches moving from the old
# Apache POI 4.1.0 output to the 5.2.3 output (true/false -> on/off, and the
# sz element relocated ahead of b/i/strike).
Set-RunFormatting "Level 2" -1 0 0 18
Set-RunFormatting "italic" 0 -1 0 0
Set-RunFormatting "bold" -1 0 0 0
Set-RunFormatting "bold & italic" -1 -1 0 0
Set-RunFormatting "Level 3" -1 0 0 14
Set-RunFormatting "Level 4" -1 0 0 12
